# Update 杭州-漫展信息 ticket-count/price data (gh-pages output regenerated at 456a3b4)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2550
$ws1.Range("F5").Value = 1444
$ws1.Range("F6").Value = 1117
$ws1.Range("F7").Value = 322
$ws1.Range("G7").Value = "已售罄"
$ws1.Range("F8").Value = 527
$ws1.Range("F11").Value = 107
$ws1.Range("F12").Value = 547
$ws1.Range("F13").Value = 8860
$ws1.Range("F18").Value = 470
$ws1.Range("F19").Value = 608
$ws1.Range("F23").Value = 2052
$ws1.Range("F24").Value = 2121
$ws1.Range("F26").Value = 1817
$ws1.Range("F28").Value = 1919
$ws1.Range("F30").Value = 329
$ws1.Range("F32").Value = 121
$ws1.Range("F34").Value = 16
$ws1.Range("F35").Value = 314
$ws1.Range("F36").Value = 60
$ws1.Range("F37").Value = 266
$ws1.Range("F38").Value = 450
$ws1.Range("F39").Value = 779
$ws1.Range("F40").Value = 272
$ws1.Range("F41").Value = 28
$ws1.Range("F43").Value = 274

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 3

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2550
$ws4.Range("F5").Value = 1444
$ws4.Range("F7").Value = 1117
$ws4.Range("F8").Value = 322
$ws4.Range("G8").Value = "已售罄"
$ws4.Range("F9").Value = 527
$ws4.Range("F11").Value = 107
$ws4.Range("F12").Value = 547
$ws4.Range("F13").Value = 8860
$ws4.Range("F19").Value = 470
$ws4.Range("F20").Value = 608
$ws4.Range("F24").Value = 2052
$ws4.Range("F25").Value = 2121
$ws4.Range("F27").Value = 1817
$ws4.Range("F29").Value = 1919
$ws4.Range("F31").Value = 329
$ws4.Range("F33").Value = 121
$ws4.Range("F35").Value = 16
$ws4.Range("F36").Value = 314
$ws4.Range("F37").Value = 60
$ws4.Range("F38").Value = 266
$ws4.Range("F39").Value = 450
$ws4.Range("F44").Value = 779
$ws4.Range("F45").Value = 3
$ws4.Range("F46").Value = 272
$ws4.Range("F47").Value = 28
$ws4.Range("F49").Value = 274
